# The row for Caso 4966 (MUÑECAS 1035) was removed from the "Optical_Power"
# sheet. All following rows (5117, 5571, 7037, 7057) shift up by one, and the
# used range shrinks from A1:N29 to A1:N28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(25).Delete()
